# Split the three "list-like" paragraphs in the Programa/Bibliografia sections
# into individual lines joined by manual line breaks (<w:br/>), matching the
# source edit that turned one run-on sentence per paragraph into one topic
# per line.
$d = $word.ActiveDocument

function Insert-LineBreaks($doc, [string[]]$items) {
    $original = [string]::Join("", $items)
    $replacement = [string]::Join("^l", $items)
    $ok = $doc.Content.Find.Execute($original, $true, $false, $false, $false, $false,
                                     $true, 1, $false, $replacement, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $($items[0])..."
    }
}

# 1. "Programa" (Portuguese) paragraph
$pt = @(
    '1. Teoria dos Sistemas. ',
    '1.1. Conceitos Básicos sobre Teoria dos Sistemas;',
    '1.2. Abordagem Sistêmica; ',
    '1.3. Classificação dos Sistemas; ',
    '1.4. Ciclo de Vida dos Sistemas; ',
    '2. Processo de Modelagem de Sistemas. ',
    '2.1. Modelo Formal; ',
    '2.2. Modelo Computacional; ',
    '2.3. Teoria dos Modelos; ',
    '3. Modelagem para Simulação ',
    '3.1. Modelagem Discreta; ',
    '3.2. Modelagem Orientada a Eventos; ',
    '4. Linguagens de Simulação. ',
    '5. Aspectos Estatísticos da Simulação de Sistemas. ',
    '5.1. Geradores de Números Aleatórios; ',
    '5.2. Geração de Variáveis Aleatórias; ',
    '5.3. Inferência Estatística; ',
    '5.4. Problemas Estatísticos Relacionados com Simulação; ',
    '6. Validação de Modelos.'
)
Insert-LineBreaks $d $pt

# 2. "Programa" (English, italic) paragraph
$en = @(
    '1. Systems Theory.',
    '1.1. Basic Concepts on Systems Theory;',
    '1.2. Systemic Approach;',
    '1.3. Classification of Systems;',
    '1.4. Systems Life Cycle;',
    '2. System Modeling Process.',
    '2.1. Formal Model;',
    '2.2. Computational Model;',
    '2.3. Theory of Models;',
    '3. Modeling for Simulation',
    '3.1. Discrete Modeling;',
    '3.2. Event-Driven Modeling;',
    '4. Simulation Languages.',
    '5. Statistical Aspects of Systems Simulation.',
    '5.1. Random Number Generators;',
    '5.2. Generation of Random Variables;',
    '5.3. Statistical inference;',
    '5.4. Statistical Problems Related to Simulation;',
    '6. Validation of Models.'
)
Insert-LineBreaks $d $en

# 3. "Bibliografia" paragraph
$bib = @(
    '1. BANKS, J., CARSON, J.S., “Discrete-Event System Simulation.”, Printice Hall, 2000.',
    '2. CHWIF, L., MEDINA, A.C. e col. “Introdução ao Simul8: um guia prático”, Livro Eletrônico, 1ª ed., 2015.',
    '3. Fishman, G.S. “Discrete-Event Simulation: Modeling, Programming, and Analysis”, Springer-Verlag, 2001.'
)
Insert-LineBreaks $d $bib

Write-Output "Programa/Bibliografia paragraphs split into lines."
